$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 757.4167
$ws.Range("J19").Value = 940.6667
$ws.Range("L19").Value = 940.6667
$ws.Range("N19").Value = -1290.6667

$ws.Range("H92").Value = 980.375
$ws.Range("I92").Value = 488
$ws.Range("K92").Value = 488
$ws.Range("M92").Value = 760

$ws.Range("H107").Value = 9570.909
$ws.Range("I107").Value = 10329
$ws.Range("K107").Value = 10329
$ws.Range("M107").Value = -8409

$ws.Range("H134").Value = 87093.664
$ws.Range("J134").Value = 89743.21000000001
$ws.Range("L134").Value = 89743.21000000001
$ws.Range("N134").Value = -99883.21000000001

$ws.Range("H135").Value = 11849.5
$ws.Range("I135").Value = 21200
$ws.Range("K135").Value = 190800
$ws.Range("M135").Value = -188265

$ws.Range("H137").Value = 476029.12
$ws.Range("I137").Value = 939031.2
$ws.Range("J137").Value = 13027.091
$ws.Range("K137").Value = 2817093.6
$ws.Range("L137").Value = 39081.273
$ws.Range("M137").Value = -2814543.6
$ws.Range("N137").Value = -44181.273

$ws.Range("H139").Value = 132185
$ws.Range("J139").Value = 132185
$ws.Range("L139").Value = 132185
$ws.Range("N139").Value = -142465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2745.9092
$ws.Range("I74").Value = 2743.7144
$ws.Range("K74").Value = 2743.7144
$ws.Range("M74").Value = -1869.7144

$ws.Range("H77").Value = 2745.9092
$ws.Range("I77").Value = 2743.7144
$ws.Range("K77").Value = 13718.572
$ws.Range("M77").Value = -9350.572

$ws.Range("H110").Value = 2745.6843
$ws.Range("I110").Value = 2434.818
$ws.Range("K110").Value = 2434.818
$ws.Range("M110").Value = -389.8180000000002

$ws.Range("H122").Value = 829355.2
$ws.Range("I122").Value = 5465.4443
$ws.Range("K122").Value = 16396.3329
$ws.Range("M122").Value = -13946.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4653.878
$ws.Range("I86").Value = 5809.9287
$ws.Range("K86").Value = 5809.9287
$ws.Range("M86").Value = -4686.9287

$ws.Range("H89").Value = 4653.878
$ws.Range("I89").Value = 5809.9287
$ws.Range("K89").Value = 29049.6435
$ws.Range("M89").Value = -23433.6435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 16668598
$ws.Range("I7").Value = 3699.6667
$ws.Range("J7").Value = 33333496
$ws.Range("K7").Value = 3699.6667
$ws.Range("L7").Value = 33333496
$ws.Range("M7").Value = -3586.6667
$ws.Range("N7").Value = -33333722

$ws.Range("H22").Value = 10989679
$ws.Range("I22").Value = 544
$ws.Range("K22").Value = 544
$ws.Range("M22").Value = -194

$ws.Range("H58").Value = 14984.909
$ws.Range("I58").Value = 17259.334
$ws.Range("K58").Value = 17259.334
$ws.Range("M58").Value = -17056.334

$ws.Range("H105").Value = 92881.44
$ws.Range("I105").Value = 124838.65
$ws.Range("J105").Value = 2336
$ws.Range("K105").Value = 124838.65
$ws.Range("L105").Value = 2336
$ws.Range("M105").Value = -123091.65
$ws.Range("N105").Value = -5830

$ws.Range("H136").Value = 14984.909
$ws.Range("I136").Value = 17259.334
$ws.Range("K136").Value = 51778.00199999999
$ws.Range("M136").Value = -49228.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5327.524
$ws.Range("I68").Value = 2124.2856
$ws.Range("J68").Value = 5968.1714
$ws.Range("K68").Value = 6372.8568
$ws.Range("L68").Value = 17904.5142
$ws.Range("M68").Value = -5561.8568
$ws.Range("N68").Value = -19526.5142

$ws.Range("H71").Value = 5327.524
$ws.Range("I71").Value = 2124.2856
$ws.Range("J71").Value = 5968.1714
$ws.Range("K71").Value = 19118.5704
$ws.Range("L71").Value = 53713.5426
$ws.Range("M71").Value = -15062.5704
$ws.Range("N71").Value = -61825.5426

$ws.Range("H113").Value = 1480.84
$ws.Range("I113").Value = 567
$ws.Range("K113").Value = 1701
$ws.Range("M113").Value = 469

$ws.Range("H131").Value = 6143.4
$ws.Range("J131").Value = 2599
$ws.Range("L131").Value = 7797
$ws.Range("N131").Value = -17877

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10050
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 10050
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 10050
$ws.Range("N10").Value = -10388
$ws.Range("M10").ClearContents()

$ws.Range("H12").Value = 3001.5
$ws.Range("I12").Value = 3001.5
$ws.Range("K12").Value = 3001.5
$ws.Range("M12").Value = -2861.5

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H55").Value = 27000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 27000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 27000
$ws.Range("N55").Value = -27654
$ws.Range("M55").ClearContents()

$ws.Range("H70").Value = 10072.357
$ws.Range("I70").Value = 9770.23
$ws.Range("J70").Value = 14000
$ws.Range("K70").Value = 9770.23
$ws.Range("L70").Value = 14000
$ws.Range("M70").Value = -9500.23
$ws.Range("N70").Value = -14540

$ws.Range("H73").Value = 10072.357
$ws.Range("I73").Value = 9770.23
$ws.Range("J73").Value = 14000
$ws.Range("K73").Value = 9770.23
$ws.Range("L73").Value = 14000
$ws.Range("M73").Value = -8834.23
$ws.Range("N73").Value = -15872

$ws.Range("H102").Value = 10258.934
$ws.Range("I102").Value = 11375.692
$ws.Range("K102").Value = 11375.692
$ws.Range("M102").Value = -9753.691999999999

$ws.Range("H113").Value = 2997.5
$ws.Range("I113").Value = 2997.5
$ws.Range("K113").Value = 2997.5
$ws.Range("M113").Value = -827.5

$ws.Range("H122").Value = 7110.6562
$ws.Range("J122").Value = 3193.25
$ws.Range("L122").Value = 9579.75
$ws.Range("N122").Value = -14479.75

$ws.Range("H132").Value = 8332.532999999999
$ws.Range("I132").Value = 8537.615
$ws.Range("K132").Value = 25612.845
$ws.Range("M132").Value = -23082.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2870.25

$ws.Range("H59").Value = 57950
$ws.Range("J59").Value = 57950
$ws.Range("L59").Value = 57950
$ws.Range("N59").Value = -59258

$ws.Range("H120").Value = 120000
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws.Range("H122").Value = 5633.1113
$ws.Range("I122").Value = 4519.6
$ws.Range("K122").Value = 13558.8
$ws.Range("M122").Value = -11108.8

$ws.Range("H136").Value = 5263.8
$ws.Range("I136").Value = 4293
$ws.Range("J136").Value = 5616.8184
$ws.Range("K136").Value = 12879
$ws.Range("L136").Value = 16850.4552
$ws.Range("M136").Value = -10329
$ws.Range("N136").Value = -21950.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14906.5
$ws.Range("J45").Value = 14906.5
$ws.Range("L45").Value = 14906.5
$ws.Range("N45").Value = -15888.5

$ws.Range("H62").Value = 91553.02
$ws.Range("I62").Value = 125267.195
$ws.Range("J62").Value = 4458.0835
$ws.Range("K62").Value = 125267.195
$ws.Range("L62").Value = 4458.0835
$ws.Range("M62").Value = -124643.195
$ws.Range("N62").Value = -5706.0835

$ws.Range("H65").Value = 91553.02
$ws.Range("I65").Value = 125267.195
$ws.Range("J65").Value = 4458.0835
$ws.Range("K65").Value = 626335.9750000001
$ws.Range("L65").Value = 22290.4175
$ws.Range("M65").Value = -623215.9750000001
$ws.Range("N65").Value = -28530.4175

$ws.Range("H107").Value = 3636.3635
$ws.Range("I107").Value = 3214.2856
$ws.Range("K107").Value = 9642.856800000001
$ws.Range("M107").Value = -7722.856800000001

$ws.Range("H122").Value = 5927.4287
$ws.Range("I122").Value = 2999.25
$ws.Range("J122").Value = 9831.666999999999
$ws.Range("K122").Value = 8997.75
$ws.Range("L122").Value = 29495.001
$ws.Range("M122").Value = -6547.75
$ws.Range("N122").Value = -34395.001

$ws.Range("H136").Value = 6383.6924
$ws.Range("I136").Value = 4887.5557
$ws.Range("K136").Value = 14662.6671
$ws.Range("M136").Value = -12112.6671
